$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(4071, 4251, 4372, 4603, 4815, 4870, 4870, 5031, 5031, 5109, 5151, 5151, 5151, 5151)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
